# Updated cryptos list with latest prices and volume deltas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.401.62"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.844.68"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.94"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6315"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07539"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2930"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.60"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07712"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "1.838.78"
$ws.Range("E12").Value = "  -7.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.001"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6799"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001046"
$ws.Range("E15").Value = "  +5.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.30"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "2.090.21"
$ws.Range("E17").Value = "  -7.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.170"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "29.424.16"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.85"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.42"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.460"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.85"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.354"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.59"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.456"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.291"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05628"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.103"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.024"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.847"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7102"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "1.247.40"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.769"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9016"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.69"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.83"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.099"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.927"
$ws.Range("E49").Value = "  -2.44%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.673"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1122"
$ws.Range("E51").Value = "  -0.47%  "
